$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 617 ("「教師要らず：10日でフランス語会話」..." entry).
# All subsequent rows (618-704) shift up by one, becoming rows 617-703,
# matching the edit recorded in the diff.
$ws.Rows("617:617").Delete()
